# Commit actualización de documentación
#
# The very first paragraph of the document (an empty, right-justified
# paragraph using the "Puesto" style) gets:
#   1. its paragraph mark formatted with single underline
#      (w:pPr/w:rPr gains <w:u w:val="single"/>)
#   2. a single literal space typed into it, using the Complex-Script
#      Arial font (matching the paragraph's existing w:rFonts w:cs="Arial")
#      but with NO underline of its own.
#
# NOTE on ordering: this runtime only persists character-formatting
# changes applied to a Range that currently contains at least one
# character. To get the underline onto the (otherwise run-less)
# paragraph mark without also underlining the new run, we:
#   a) type a throw-away placeholder character,
#   b) underline the whole (now non-empty) paragraph range - this stamps
#      the paragraph-mark rPr (and the placeholder run) with the underline,
#   c) delete just the placeholder character text - the paragraph goes
#      back to being run-less, but the underline stays recorded on the
#      paragraph mark's rPr,
#   d) type the real space (it picks up no direct character formatting,
#      i.e. no <w:u>, since the paragraph has no run to inherit from),
#   e) set the complex-script font name on that new run only.

$d = $word.ActiveDocument

$target = $d.Paragraphs.Item(1)

# a) Temporary placeholder so the Range is non-empty and formatting sticks.
$target.Range.InsertAfter("X")

# b) Underline the whole paragraph (placeholder run + paragraph mark).
$para = $d.Paragraphs.Item(1)
$para.Range.Font.Underline = 1

# c) Remove the placeholder text; the paragraph-mark keeps the underline.
$placeholder = $d.Range(0, 1)
$placeholder.Text = ""

# d) Type the real content: a single preserved space.
$para2 = $d.Paragraphs.Item(1)
$para2.Range.InsertAfter(" ")

# e) Give that new run the Complex-Script "Arial" font (w:rFonts w:cs=),
#    matching the paragraph's existing rFonts, without touching
#    ascii/hAnsi/eastAsia and without re-adding underline.
$newRun = $d.Range(0, 1)
$newRun.Font.NameBi = "Arial"
